# Update gh-pages to output generated at 456a3b4
# Apply updated "want to go" counts (F) / minimum price (G) figures to the
# convention-info workbook. The 4th sheet ("全部类型") aggregates rows from
# the other sheets, so matching rows are updated in both places.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet "展览" (sheet1) ---
$ws1.Range("F2").Value = 647
$ws1.Range("G2").Value = 55

$ws1.Range("G3").Value = "不可售"

$ws1.Range("F5").Value = 591

$ws1.Range("F7").Value = 3259

$ws1.Range("F8").Value = 479

$ws1.Range("F9").Value = 8304

$ws1.Range("F11").Value = 483

$ws1.Range("F13").Value = 470

# --- Sheet "演出" (sheet2) ---
$ws2.Range("F4").Value = 3

# --- Sheet "全部类型" (sheet4, aggregates all rows) ---
$ws4.Range("F2").Value = 647
$ws4.Range("G2").Value = 55

$ws4.Range("G3").Value = "不可售"

$ws4.Range("F5").Value = 591

$ws4.Range("F9").Value = 3259

$ws4.Range("F10").Value = 479

$ws4.Range("F11").Value = 3

$ws4.Range("F12").Value = 8304

$ws4.Range("F14").Value = 483

$ws4.Range("F18").Value = 470
